# The "<id>p017r_1</id>" tag is split across three runs in the source doc:
#   <id>      (Courier New, color 7f6000, size 9pt)
#   p017r_1   (default run formatting, color 000000)
#   </id>     (Courier New, color 7f6000, size 9pt)
# The edit collapses them into a single run containing the concatenated
# text "<id>p017r_1</id>", keeping the formatting of the first ("<id>")
# run - exactly like Word does when you assign Range.Text for a range
# that spans several runs (the new text adopts the formatting of the
# range's first run, and the now-empty extra runs are dropped).

$d = $word.ActiveDocument

$target = "<id>p017r_1</id>"

$rng = $d.Content
$found = $rng.Find.Execute($target, $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)

if ($found) {
    # Re-assigning Range.Text with the exact same string the range already
    # holds is a no-op, so the runs would stay split. Nudge the text to a
    # scratch value first, then set it to the real target so the range
    # actually gets rewritten as a single merged run.
    $rng.Text = "##SCRATCH##"

    $rng2 = $d.Content
    [void]$rng2.Find.Execute("##SCRATCH##", $true, $false, $false, $false, `
                              $false, $true, 1, $false, "", 0)
    $rng2.Text = $target
} else {
    throw "Could not find target text '$target' to merge runs."
}
